$d = $word.ActiveDocument

$replacements = @(
    @("55×49=", "25×91="),
    @("29×21=", "66×60="),
    @("46×16=", "73×32="),
    @("95×15=", "45×70="),
    @("86×37=", "12×32="),
    @("57×38=", "64×52="),
    @("26×68=", "56×34="),
    @("98×64=", "55×58="),
    @("84×94=", "80×15="),
    @("33×76=", "41×69="),
    @("53×35=", "41×64="),
    @("35×14=", "60×49="),
    @("30×71=", "86×88="),
    @("74×60=", "38×40="),
    @("73×96=", "19×25="),
    @("40×69=", "51×22="),
    @("19×90=", "69×21="),
    @("92×72=", "55×96="),
    @("19×45=", "88×35="),
    @("92×45=", "74×71="),
    @("34×13=", "43×62="),
    @("89×35=", "87×16="),
    @("25×78=", "58×47="),
    @("42×15=", "29×52="),
    @("50×57=", "26×74=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
